$wb = $excel.ActiveWorkbook

# --- Add Sheet2 and Sheet3 after Sheet1 (so order is Sheet1, Sheet2, Sheet3) ---
$sheet1 = $wb.Worksheets.Item(1)
$sheet2 = $wb.Worksheets.Add($null, $sheet1)
$sheet3 = $wb.Worksheets.Add($null, $sheet2)

# =========================================================================
# Sheet2: demonstrates the various number-format categories
# =========================================================================
$sheet2.Range("A1").Value = "General"
$sheet2.Range("B1").Value = 12345

$sheet2.Range("A2").Value = "Number"
$sheet2.Range("B2").Value = 12345
$sheet2.Range("B2").NumberFormat = "0.00"

$sheet2.Range("A3").Value = "Currency"
$sheet2.Range("B3").Value = 12345
$sheet2.Range("B3").NumberFormat = "#,##0.00\ ""₽"""

$sheet2.Range("A4").Value = "Accounting"
$sheet2.Range("B4").Value = 12345
$sheet2.Range("B4").NumberFormat = "_-* #,##0.00\ ""₽""_-;\-* #,##0.00\ ""₽""_-;_-* ""-""??\ ""₽""_-;_-@_-"

$sheet2.Range("A5").Value = "Date"
$sheet2.Range("B5").Value = 12345
$sheet2.Range("B5").NumberFormat = "m/d/yy"

$sheet2.Range("A6").Value = "Time"
$sheet2.Range("B6").Value = 12345
$sheet2.Range("B6").NumberFormat = "[$-F400]h:mm:ss\ AM/PM"

$sheet2.Range("A7").Value = "Percentage"
$sheet2.Range("B7").Value = 12345
$sheet2.Range("B7").NumberFormat = "0.00%"

$sheet2.Range("A8").Value = "Fraction"
$sheet2.Range("B8").Value = 12345
$sheet2.Range("B8").NumberFormat = "# ?/?"

$sheet2.Range("A9").Value = "Scientific"
$sheet2.Range("B9").Value = 12345
$sheet2.Range("B9").NumberFormat = "0.00E+00"

$sheet2.Range("A10").Value = "Text"
$sheet2.Range("B10").Value = 12345
$sheet2.Range("B10").NumberFormat = "@"

$sheet2.Columns.Item(1).ColumnWidth = 10.05078125
$sheet2.Columns.Item(2).ColumnWidth = 11.5234375

$sheet2.Range("B10").Select()

# =========================================================================
# Sheet3: demonstrates the "text cell should not be converted to a number" fix
# =========================================================================
$sheet3.Range("A1").Value = 45325
$sheet3.Range("A1").NumberFormat = "d/m/yy;@"

$sheet3.Range("A2").Value = 45325
$sheet3.Range("A2").NumberFormat = "m/d/yy"

$sheet3.Range("A3").Value = "3"
$sheet3.Range("A3").NumberFormat = "m/d/yy"
$sheet3.Range("A4").Value = "3.2"
$sheet3.Range("A4").NumberFormat = "m/d/yy"
$sheet3.Range("A5").Value = "3.2.24"
$sheet3.Range("A5").NumberFormat = "m/d/yy"
$sheet3.Range("A6").Value = "3.2.24.7"
$sheet3.Range("A6").NumberFormat = "m/d/yy"
$sheet3.Range("A7").Value = "3.2.24.d"
$sheet3.Range("A7").NumberFormat = "m/d/yy"

$sheet3.Columns.Item(1).ColumnWidth = 9.62890625

# Sheet3 is the active sheet (it was the last one added); select its first cell
$sheet3.Range("A1").Select()
$sheet3.Activate()
